# Apply the NOTIFY_SEND_SETTING changes described by the commit:
#  - add a new "IS_ALLOW_CONFIG" column (M) to NOTIFY_SEND_SETTING, with
#    values for each existing data row.
#  - update the FD_LEVEL column (I) so some rows that used to say "site"
#    now correctly say "organization" or "project" (addUser/enableOrganization/
#    disableOrganization -> organization; enableProject/disableProject -> project),
#    matching the "remove send setting organization and modify send setting
#    site" part of the commit message.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("NOTIFY_SEND_SETTING")

# New header for the added column
$ws2.Range("M7").Value = "IS_ALLOW_CONFIG"

# FD_LEVEL (column I) corrections
$ws2.Range("I9").Value  = "organization"   # addUser
$ws2.Range("I12").Value = "project"        # enableProject
$ws2.Range("I13").Value = "project"        # disableProject
$ws2.Range("I14").Value = "organization"   # enableOrganization
$ws2.Range("I15").Value = "organization"   # disableOrganization

# IS_ALLOW_CONFIG (column M) values for each data row
$ws2.Range("M8").Value  = 0   # forgetPassword
$ws2.Range("M9").Value  = 1   # addUser
$ws2.Range("M10").Value = 1   # addFunction
$ws2.Range("M11").Value = 1   # modifyPassword
$ws2.Range("M12").Value = 1   # enableProject
$ws2.Range("M13").Value = 1   # disableProject
$ws2.Range("M14").Value = 1   # enableOrganization
$ws2.Range("M15").Value = 1   # disableOrganization
$ws2.Range("M16").Value = 0   # registerOrganization
